$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '24.593.59'
Set-TextValue $ws.Range('E2') '  -0.32%  '
Set-TextValue $ws.Range('D3') '1.689.43'
Set-TextValue $ws.Range('E3') '  +0.01%  '
Set-TextValue $ws.Range('E4') '  +0.10%  '
Set-TextValue $ws.Range('D5') '313.66'
Set-TextValue $ws.Range('E5') '  -0.53%  '
Set-TextValue $ws.Range('E6') '  +0.12%  '
Set-TextValue $ws.Range('D7') '0.3899'
Set-TextValue $ws.Range('E7') '  -1.03%  '
Set-TextValue $ws.Range('D8') '0.4032'
Set-TextValue $ws.Range('E8') '  -0.40%  '
Set-TextValue $ws.Range('D9') '1.496'
Set-TextValue $ws.Range('E9') '  +0.65%  '
Set-TextValue $ws.Range('D10') '1.004'
Set-TextValue $ws.Range('E10') '  +0.18%  '
Set-TextValue $ws.Range('D11') '52.61'
Set-TextValue $ws.Range('E11') '  +0.11%  '
Set-TextValue $ws.Range('D12') '0.08750'
Set-TextValue $ws.Range('E12') '  -0.86%  '
Set-TextValue $ws.Range('D13') '7.586'
Set-TextValue $ws.Range('E13') '  +4.86%  '
Set-TextValue $ws.Range('D14') '24.82'
Set-TextValue $ws.Range('E14') '  +5.74%  '
Set-TextValue $ws.Range('D15') '7.957'
Set-TextValue $ws.Range('E15') '  -0.94%  '
Set-TextValue $ws.Range('E16') '  +2.55%  '
Set-TextValue $ws.Range('D17') '1.681.57'
Set-TextValue $ws.Range('E17') '  -0.52%  '
Set-TextValue $ws.Range('D18') '98.47'
Set-TextValue $ws.Range('E18') '  -1.03%  '
Set-TextValue $ws.Range('D19') '0.07107'
Set-TextValue $ws.Range('E19') '  +1.40%  '
Set-TextValue $ws.Range('D20') '19.84'
Set-TextValue $ws.Range('E20') '  +1.81%  '
Set-TextValue $ws.Range('D21') '7.277'
Set-TextValue $ws.Range('E21') '  +4.21%  '
Set-TextValue $ws.Range('E22') '  -0.30%  '
Set-TextValue $ws.Range('D23') '14.23'
Set-TextValue $ws.Range('E23') '  -0.45%  '
Set-TextValue $ws.Range('D24') '24.593.15'
Set-TextValue $ws.Range('E24') '  -0.22%  '
Set-TextValue $ws.Range('D25') '3.000'
Set-TextValue $ws.Range('E25') '  -9.15%  '
Set-TextValue $ws.Range('D26') '2.355'
Set-TextValue $ws.Range('E26') '  -0.02%  '
Set-TextValue $ws.Range('D27') '22.79'
Set-TextValue $ws.Range('E27') '  +0.24%  '
Set-TextValue $ws.Range('D28') '161.80'
Set-TextValue $ws.Range('E28') '  -0.45%  '
Set-TextValue $ws.Range('D29') '8.749'
Set-TextValue $ws.Range('E29') '  +14.70%  '
Set-TextValue $ws.Range('D30') '136.58'
Set-TextValue $ws.Range('E30') '  +0.91%  '
Set-TextValue $ws.Range('D31') '5.221'
Set-TextValue $ws.Range('E31') '  +0.85%  '
Set-TextValue $ws.Range('D32') '1.866.36'
Set-TextValue $ws.Range('E32') '  -0.79%  '
Set-TextValue $ws.Range('D33') '0.08802'
Set-TextValue $ws.Range('E33') '  +3.16%  '
Set-TextValue $ws.Range('D34') '7.375'
Set-TextValue $ws.Range('E34') '  +3.96%  '
Set-TextValue $ws.Range('D35') '1.038'
Set-TextValue $ws.Range('E35') '  -1.82%  '
Set-TextValue $ws.Range('D36') '1.990'
Set-TextValue $ws.Range('E36') '  +5.42%  '
Set-TextValue $ws.Range('D37') '0.02916'
Set-TextValue $ws.Range('E37') '  +7.24%  '
Set-TextValue $ws.Range('D38') '0.2727'
Set-TextValue $ws.Range('E38') '  -0.17%  '
Set-TextValue $ws.Range('D39') '10.75'
Set-TextValue $ws.Range('E39') '  -3.98%  '
Set-TextValue $ws.Range('B40') 'Stellar'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D40') '0.09137'
Set-TextValue $ws.Range('E40') '  -0.45%  '
Set-TextValue $ws.Range('B41') 'Aptos'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D41') '14.19'
Set-TextValue $ws.Range('E41') '  -1.63%  '
Set-TextValue $ws.Range('D42') '0.7804'
Set-TextValue $ws.Range('E42') '  +2.44%  '
Set-TextValue $ws.Range('D43') '1.455'
Set-TextValue $ws.Range('E43') '  -0.64%  '
Set-TextValue $ws.Range('D44') '16.61'
Set-TextValue $ws.Range('E44') '  +4.10%  '
Set-TextValue $ws.Range('D45') '0.7191'
Set-TextValue $ws.Range('E45') '  +0.91%  '
Set-TextValue $ws.Range('D46') '2.585'
Set-TextValue $ws.Range('E46') '  +0.09%  '
Set-TextValue $ws.Range('D47') '4.189'
Set-TextValue $ws.Range('E47') '  -0.49%  '
Set-TextValue $ws.Range('E48') '  +0.11%  '
Set-TextValue $ws.Range('D49') '1.329'
Set-TextValue $ws.Range('E49') '  +0.95%  '
Set-TextValue $ws.Range('D50') '137.63'
Set-TextValue $ws.Range('E50') '  -1.56%  '
Set-TextValue $ws.Range('D51') '90.90'
Set-TextValue $ws.Range('E51') '  +1.35%  '
